$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.090.31"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.45%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.521.51"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.72%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.27"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.80%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.521.51"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.14%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("E9").Value = "  -1.47%  "

# Row 10
$ws.Range("E10").Value = "  -0.46%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.94"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.78%  "

# Row 12
$ws.Range("E12").Value = "  -1.58%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000218"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.77%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.08"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.08%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.116.16"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.96%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.518.59"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.062.81"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +9.78%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.34"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.611"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.51%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.55%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.660.77"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.84%  "

# Row 26
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.80%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.35%  "

# Row 30
$ws.Range("E30").Value = "  +0.57%  "

# Row 31
$ws.Range("E31").Value = "  -2.82%  "

# Row 32
$ws.Range("E32").Value = "  -2.25%  "

# Row 33
$ws.Range("E33").Value = "  +0.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.49"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.512.10"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.91%  "

# Row 36
$ws.Range("E36").Value = "  -2.34%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.95"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.68%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.05"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.22%  "

# Row 39
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "174.45"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.18%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0896"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -9.42%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.897"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.03"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.94%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.29"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.63%  "

# Row 49
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.47"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.58%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.49"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.992"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.72%  "
